# Update the grants/sub-awards table on Sheet3 with new placeholder data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Column D switches from a text grant-ID (shared string) to a plain numeric id.
$ws.Range("D2").Value = 100
$ws.Range("D3").Value = 101
$ws.Range("D4").Value = 102

# Column E: PI name -> generic placeholder name.
$ws.Range("E2").Value = "name 1"
$ws.Range("E3").Value = "name 2"
$ws.Range("E4").Value = "name 3"

# Column F: grant title -> generic placeholder title.
$ws.Range("F2").Value = "Title 1"
$ws.Range("F3").Value = "Title 2 "
$ws.Range("F4").Value = "Title 3"

# Move / refresh the active selection to F5, as in the saved workbook.
$ws.Activate()
$ws.Range("F5").Select()
